$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates for column F ("想去人数")
$sheetUpdates = @{
    "展览" = @(
        @{Row = 2;  Value = 1075},
        @{Row = 3;  Value = 359},
        @{Row = 4;  Value = 1470},
        @{Row = 5;  Value = 8691},
        @{Row = 7;  Value = 488},
        @{Row = 8;  Value = 638},
        @{Row = 9;  Value = 279},
        @{Row = 12; Value = 3537},
        @{Row = 13; Value = 48},
        @{Row = 15; Value = 75},
        @{Row = 16; Value = 1134},
        @{Row = 20; Value = 196},
        @{Row = 21; Value = 2266},
        @{Row = 22; Value = 50}
    )
    "全部类型" = @(
        @{Row = 2;  Value = 1075},
        @{Row = 3;  Value = 359},
        @{Row = 4;  Value = 1470},
        @{Row = 5;  Value = 8691},
        @{Row = 7;  Value = 488},
        @{Row = 8;  Value = 638},
        @{Row = 9;  Value = 279},
        @{Row = 12; Value = 3537},
        @{Row = 13; Value = 48},
        @{Row = 15; Value = 75},
        @{Row = 16; Value = 1134},
        @{Row = 20; Value = 196},
        @{Row = 21; Value = 2266},
        @{Row = 23; Value = 50}
    )
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($update in $sheetUpdates[$sheetName]) {
        $ws.Cells.Item($update.Row, 6).Value = $update.Value
    }
}
